$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = "2024-09-25T18:06:40Z"
$ws.Range("B51").Value = "temperature"

# "25" reads as a number to Excel's type-inference; force text storage
# (matching the rest of the column) then drop back to the Normal style
# so no stray formatting is left behind on the cell.
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "25"
$ws.Range("C51").Style = "Normal"

$ws.Range("D51").Value = "N/A"
$ws.Range("E51").Value = "N/A"
$ws.Range("F51").Value = "N/A"
